# Generate Report for Handoff
# Updates the latest handoff/handback timestamps for the
# f3a64252-2f7d-4f38-8027-b2e97acbdb85 row (row 7) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-37-13 02:37:30"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-13 02:37:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-13 02:37:30"
